$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Zoom level: 115% -> 70% ---
$excel.ActiveWindow.Zoom = 70

# --- New header row values (R3:U3 = 2,3,4,5) ---
$ws.Range("R3").Value2 = 2
$ws.Range("S3").Value2 = 3
$ws.Range("T3").Value2 = 4
$ws.Range("U3").Value2 = 5

# --- New frequency-count formulas filled the same way a user would:
#     type R4 alone, fill S4:U19 as one block, then R5:U31 as another block.
#     (Mirrors the exact shared-formula grouping Excel itself produces.) ---
$ws.Range("R4").Formula = "=IF(R`$3=`$M4,1,0)"
$ws.Range("S4:U19").Formula = "=IF(S`$3=`$M4,1,0)"
$ws.Range("R5:U31").Formula = "=IF(R`$3=`$M5,1,0)"

# --- New color-scale conditional formatting over R4:U31, pushed to top priority ---
$cs = $ws.Range("R4:U31").FormatConditions.AddColorScale(3)
$cs.SetFirstPriority()

# --- Move the active selection to W11 ---
$ws.Range("W11").Select() | Out-Null
